$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12
$ws.Range("C3").Value = 11
$ws.Range("C5").Value = 20

$ws.Range("C2").Select()
